$d = $word.ActiveDocument

$map = @{
    "MR.Bella Abdelouahab" = "MR.bella abdelouahab"
    "trdnt N 28 Agadir"    = "traoudant no where"
    "787898"               = "+212762549778"
    "       RR454545454MA" = "       RR000000023MA"
    "MR.Bella Maha"        = "MR.nas ons"
    "trdnt N 29 Agadir"    = "5dsd"
    "654987"               = "5684531"
}

foreach ($p in $d.Paragraphs) {
    $r = $p.Range
    $t = $r.Text.TrimEnd([char]13, [char]7)
    if ($map.ContainsKey($t)) {
        $r.Text = $map[$t]
    }
}
